$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1226.3889
$ws.Range("I98").Value = 1226.3889
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1226.3889
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 271.6111000000001

# Row 107
$ws.Range("H107").Value = 500370.84
$ws.Range("I107").Value = 625306.2
$ws.Range("J107").Value = 629.5
$ws.Range("K107").Value = 625306.2
$ws.Range("L107").Value = 629.5
$ws.Range("M107").Value = -623386.2
$ws.Range("N107").Value = -4469.5

# Row 122
$ws.Range("H122").Value = 1226.3889
$ws.Range("I122").Value = 1226.3889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3679.1667
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1229.1667

# Row 132
$ws.Range("H132").Value = 1773.0588
$ws.Range("I132").Value = 1379.3096
$ws.Range("K132").Value = 4137.9288
$ws.Range("M132").Value = -1607.9288

# Row 137
$ws.Range("H137").Value = 9617163
$ws.Range("I137").Value = 1556.091
$ws.Range("J137").Value = 26317952
$ws.Range("K137").Value = 4668.272999999999
$ws.Range("L137").Value = 78953856
$ws.Range("M137").Value = -2118.272999999999
$ws.Range("N137").Value = -78958956

# Row 138
$ws.Range("H138").Value = 2770.48
$ws.Range("I138").Value = 1112.7646
$ws.Range("J138").Value = 3624.4546
$ws.Range("K138").Value = 3338.2938
$ws.Range("L138").Value = 10873.3638
$ws.Range("M138").Value = 1801.7062
$ws.Range("N138").Value = -21153.3638

# Row 141
$ws.Range("H141").Value = 1570
$ws.Range("I141").Value = 1570
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4710
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 470

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 39345.223
$ws.Range("I2").Value = 47114.59
$ws.Range("K2").Value = 47114.59
$ws.Range("M2").Value = -47001.59

# Row 32
$ws.Range("H32").Value = 24620.627
$ws.Range("I32").Value = 22247.344
$ws.Range("K32").Value = 22247.344
$ws.Range("M32").Value = -21960.344

# Row 45
$ws.Range("H45").Value = 2538.389
$ws.Range("I45").Value = 1499.1818
$ws.Range("J45").Value = 4171.4287
$ws.Range("K45").Value = 1499.1818
$ws.Range("L45").Value = 4171.4287
$ws.Range("M45").Value = -1122.1818
$ws.Range("N45").Value = -4925.4287

# Row 61
$ws.Range("H61").Value = 1571.8918
$ws.Range("I61").Value = 1349.28
$ws.Range("J61").Value = 2035.6666
$ws.Range("K61").Value = 1349.28
$ws.Range("L61").Value = 2035.6666
$ws.Range("M61").Value = -1137.28
$ws.Range("N61").Value = -2459.6666

# Row 116
$ws.Range("H116").Value = 39345.223
$ws.Range("I116").Value = 47114.59
$ws.Range("K116").Value = 47114.59
$ws.Range("M116").Value = -44820.59

# Row 122
$ws.Range("H122").Value = 1183.6
$ws.Range("I122").Value = 1183.6
$ws.Range("K122").Value = 3550.8
$ws.Range("M122").Value = -1100.8

# Row 136
$ws.Range("H136").Value = 1571.8918
$ws.Range("I136").Value = 1349.28
$ws.Range("J136").Value = 2035.6666
$ws.Range("K136").Value = 4047.84
$ws.Range("L136").Value = 6106.9998
$ws.Range("M136").Value = -1497.84
$ws.Range("N136").Value = -11206.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 39345.223
$ws.Range("I3").Value = 47114.59
$ws.Range("K3").Value = 47114.59
$ws.Range("M3").Value = -47000.59

# Row 53
$ws.Range("H53").Value = 39333.332
$ws.Range("J53").Value = 39333.332
$ws.Range("L53").Value = 39333.332
$ws.Range("N53").Value = -40481.332

# Row 86
$ws.Range("H86").Value = 2024.2333
$ws.Range("I86").Value = 1836.3158
$ws.Range("J86").Value = 2348.818
$ws.Range("K86").Value = 1836.3158
$ws.Range("L86").Value = 2348.818
$ws.Range("M86").Value = -713.3158000000001
$ws.Range("N86").Value = -4594.818

# Row 89
$ws.Range("H89").Value = 2024.2333
$ws.Range("I89").Value = 1836.3158
$ws.Range("J89").Value = 2348.818
$ws.Range("K89").Value = 9181.579
$ws.Range("L89").Value = 11744.09
$ws.Range("M89").Value = -3565.579
$ws.Range("N89").Value = -22976.09

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 344.4375
$ws.Range("I107").Value = 274.81818
$ws.Range("J107").Value = 380.90475
$ws.Range("K107").Value = 274.81818
$ws.Range("L107").Value = 380.90475
$ws.Range("M107").Value = 1645.18182
$ws.Range("N107").Value = -4220.90475

# Row 132
$ws.Range("H132").Value = 1801.1177
$ws.Range("I132").Value = 1062.5385
$ws.Range("J132").Value = 4201.5
$ws.Range("K132").Value = 3187.6155
$ws.Range("L132").Value = 12604.5
$ws.Range("M132").Value = -657.6155000000003
$ws.Range("N132").Value = -17664.5

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 439.23077
$ws.Range("I122").Value = 376.86365
$ws.Range("J122").Value = 782.25
$ws.Range("K122").Value = 3391.77285
$ws.Range("L122").Value = 7040.25
$ws.Range("M122").Value = -941.7728500000003
$ws.Range("N122").Value = -11940.25

# Row 131
$ws.Range("H131").Value = 11930731
$ws.Range("I131").Value = 29471560
$ws.Range("J131").Value = 2967.64
$ws.Range("K131").Value = 88414680
$ws.Range("L131").Value = 8902.92
$ws.Range("M131").Value = -88409640
$ws.Range("N131").Value = -18982.92

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1445255.8
$ws.Range("I122").Value = 1857471.8
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5572415.4
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -5569965.4
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 5168.643
$ws.Range("I122").Value = 5655.522
$ws.Range("J122").Value = 2929
$ws.Range("K122").Value = 16966.566
$ws.Range("L122").Value = 8787
$ws.Range("M122").Value = -14516.566
$ws.Range("N122").Value = -13687

# Row 132
$ws.Range("H132").Value = 1490280.5
$ws.Range("I132").Value = 1940025.8
$ws.Range("J132").Value = 2661.7693
$ws.Range("K132").Value = 5820077.4
$ws.Range("L132").Value = 7985.3079
$ws.Range("M132").Value = -5817547.4
$ws.Range("N132").Value = -13045.3079

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2632.4468
$ws.Range("I136").Value = 2882.4138
$ws.Range("J136").Value = 2229.7222
$ws.Range("K136").Value = 8647.241399999999
$ws.Range("L136").Value = 6689.1666
$ws.Range("M136").Value = -6097.241399999999
$ws.Range("N136").Value = -11789.1666

# Row 139
$ws.Range("H139").Value = 172875
$ws.Range("J139").Value = 172875
$ws.Range("L139").Value = 172875
$ws.Range("N139").Value = -183155
